$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 9576478
$ws.Range("B2").Value = 42615
$ws.Range("C2").Value = "Vitaliy provided very informative directions and instructions on getting into the apartment. The apartment was exactly like it was in the photos. Harlem is a unique neighbourhood in NYC and will give you a different flavour of the city. "

# Update row 3
$ws.Range("A3").Value = 27930717
$ws.Range("B3").Value = 43347
$ws.Range("C3").Value = "The location of this apartment is great - very close to the subway and it allows you to reach Manhattan quickly. Good restaurants in the neighborhood. Maxime's (Hidden by Airbnb) , who were at home when we stayed, were absolutely fantastic and they made our experience in NYC special."

# Update row 4
$ws.Range("A4").Value = 13192097
$ws.Range("B4").Value = 43524
$ws.Range("C4").Value = "A rather cool yet inviting & comfortable place.  We felt like we had our own little nest in the big city!❤️"

# Update row 5
$ws.Range("A5").Value = 18803064
$ws.Range("B5").Value = 43411
$ws.Range("C5").Value = "Einfach nur empfehlenswert!!!"
